$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update changed values
$ws.Range("D2").Value = 10908
$ws.Range("E2").Value = -724
$ws.Range("F2").Value = -724
$ws.Range("G2").Value = -1378
$ws.Range("H2").Value = -1069
$ws.Range("I2").Value = -1050
$ws.Range("J2").Value = -19
$ws.Range("K2").Value = 13501
$ws.Range("L2").Value = 11597
$ws.Range("M2").Value = 1905
$ws.Range("N2").Value = 1859
$ws.Range("O2").Value = 45
$ws.Range("P2").Value = 496
$ws.Range("Q2").Value = 673
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = -836
$ws.Range("T2").Value = 4
$ws.Range("U2").Value = 669
$ws.Range("V2").Value = 3832
$ws.Range("W2").Value = -6.64
$ws.Range("X2").Value = -9.800000000000001
$ws.Range("Y2").Value = -44.52
$ws.Range("Z2").Value = -7.38
$ws.Range("AA2").Value = 608.88
$ws.Range("AB2").Value = 261.37
$ws.Range("AC2").Value = -10589
$ws.Range("AD2").Value = -1.19
$ws.Range("AE2").Value = 18759
$ws.Range("AF2").Value = 0.67
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 9908391

# Row 3: update changed values
$ws.Range("D3").Value = 13581
$ws.Range("E3").Value = 409
$ws.Range("F3").Value = 409
$ws.Range("G3").Value = 318
$ws.Range("H3").Value = 222
$ws.Range("I3").Value = 206
$ws.Range("J3").Value = 16
$ws.Range("K3").Value = 15294
$ws.Range("L3").Value = 12960
$ws.Range("M3").Value = 2335
$ws.Range("N3").Value = 2284
$ws.Range("O3").Value = 51
$ws.Range("P3").Value = 496
$ws.Range("Q3").Value = 2640
$ws.Range("R3").Value = 35
$ws.Range("S3").Value = -2221
$ws.Range("T3").Value = 7
$ws.Range("U3").Value = 2633
$ws.Range("V3").Value = 1901
$ws.Range("W3").Value = 3.01
$ws.Range("X3").Value = 1.63
$ws.Range("Y3").Value = 9.93
$ws.Range("Z3").Value = 1.54
$ws.Range("AA3").Value = 555.05
$ws.Range("AB3").Value = 293.68
$ws.Range("AC3").Value = 2076
$ws.Range("AD3").Value = 9.06
$ws.Range("AE3").Value = 23042
$ws.Range("AF3").Value = 0.82
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 9908391

# Row 4: update changed values
$ws.Range("D4").Value = 17723
$ws.Range("E4").Value = 697
$ws.Range("F4").Value = 697
$ws.Range("G4").Value = 448
$ws.Range("H4").Value = 267
$ws.Range("I4").Value = 265
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 15534
$ws.Range("L4").Value = 12853
$ws.Range("M4").Value = 2681
$ws.Range("N4").Value = 2629
$ws.Range("O4").Value = 53
$ws.Range("P4").Value = 508
$ws.Range("Q4").Value = 737
$ws.Range("R4").Value = -427
$ws.Range("S4").Value = 163
$ws.Range("T4").Value = 9
$ws.Range("U4").Value = 728
$ws.Range("V4").Value = 2262
$ws.Range("W4").Value = 3.93
$ws.Range("X4").Value = 1.5
$ws.Range("Y4").Value = 10.8
$ws.Range("Z4").Value = 1.73
$ws.Range("AA4").Value = 479.38
$ws.Range("AB4").Value = 330.3
$ws.Range("AC4").Value = 2652
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 25847
$ws.Range("AF4").Value = 0.62
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 1.57
$ws.Range("AI4").Value = 9.59
$ws.Range("AJ4").Value = 10166455

# Row 5: update changed values
$ws.Range("D5").Value = 19843
$ws.Range("E5").Value = 1333
$ws.Range("F5").Value = 1333
$ws.Range("G5").Value = 977
$ws.Range("H5").Value = 725
$ws.Range("I5").Value = 720
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 19054
$ws.Range("L5").Value = 15541
$ws.Range("M5").Value = 3513
$ws.Range("N5").Value = 3460
$ws.Range("O5").Value = 52
$ws.Range("P5").Value = 545
$ws.Range("Q5").Value = -2028
$ws.Range("R5").Value = 35
$ws.Range("S5").Value = 2326
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = -2038
$ws.Range("V5").Value = 4748
$ws.Range("W5").Value = 6.72
$ws.Range("X5").Value = 3.65
$ws.Range("Y5").Value = 23.64
$ws.Range("Z5").Value = 4.19
$ws.Range("AA5").Value = 442.41
$ws.Range("AB5").Value = 455.25
$ws.Range("AC5").Value = 6840
$ws.Range("AD5").Value = 3.32
$ws.Range("AE5").Value = 31750
$ws.Range("AF5").Value = 0.71
$ws.Range("AG5").Value = 375
$ws.Range("AH5").Value = 1.65
$ws.Range("AI5").Value = 5.68
$ws.Range("AJ5").Value = 10895933

# Row 6: update changed values
$ws.Range("D6").Value = 21422
$ws.Range("E6").Value = 2145
$ws.Range("F6").Value = 2145
$ws.Range("G6").Value = 2271
$ws.Range("H6").Value = 1620
$ws.Range("I6").Value = 1621
$ws.Range("K6").Value = 17784
$ws.Range("L6").Value = 13324
$ws.Range("M6").Value = 4460
$ws.Range("N6").Value = 4393
$ws.Range("P6").Value = 579
$ws.Range("Q6").Value = 137
$ws.Range("R6").Value = -363
$ws.Range("S6").Value = -526
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 127
$ws.Range("V6").Value = 4370
$ws.Range("W6").Value = 10.01
$ws.Range("X6").Value = 7.56
$ws.Range("Y6").Value = 41.29
$ws.Range("Z6").Value = 8.800000000000001
$ws.Range("AA6").Value = 298.72
$ws.Range("AB6").Value = 600.21
$ws.Range("AC6").Value = 14139
$ws.Range("AD6").Value = 1.19
$ws.Range("AE6").Value = 37956
$ws.Range("AF6").Value = 0.44
$ws.Range("AG6").Value = 375
$ws.Range("AH6").Value = 2.23
$ws.Range("AI6").Value = 2.68
$ws.Range("AJ6").Value = 11570702

# Row 7: update changed values
$ws.Range("D7").Value = 15737
$ws.Range("E7").Value = 872
$ws.Range("G7").Value = 797
$ws.Range("I7").Value = 700
$ws.Range("W7").Value = 5.54
$ws.Range("AC7").Value = 6048
$ws.Range("AD7").Value = 2.24

# Row 7: clear removed cells
$ws.Range("H7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: update changed values
$ws.Range("D8").Value = 14910
$ws.Range("E8").Value = 740
$ws.Range("G8").Value = 728
$ws.Range("I8").Value = 553
$ws.Range("W8").Value = 4.96
$ws.Range("AC8").Value = 4778
$ws.Range("AD8").Value = 2.84

# Row 8: clear removed cells
$ws.Range("H8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: clear removed cells
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
